# The deck ships two theme parts:
#   theme1.xml - "Office Theme" colours (currently only wired to the Notes
#                Master, so it is not reachable through the Slide/Notes
#                Master object model exposed here)
#   theme2.xml - "Integral" colours (wired to the one-and-only Slide Master,
#                i.e. this is the theme that actually paints every slide and
#                is the one exposed by SlideMaster.Theme)
#
# The authored change swaps which named colour set lives in which theme
# part so the deck's visible design switches from "Integral" back to the
# default "Office Theme" palette. We reproduce that visible effect by
# repointing the Slide Master's live ThemeColorScheme to the twelve
# standard Office Theme colours, in the documented
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink slot order.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      (hex 44546A, byte-swapped for RGB long)
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      (hex E7E6E6)
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  (hex 5B9BD5)
$tcs.Item(6).RGB  = 0x317DED   # accent2  (hex ED7D31)
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  (hex A5A5A5)
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  (hex FFC000)
$tcs.Item(9).RGB  = 0xC47244   # accent5  (hex 4472C4)
$tcs.Item(10).RGB = 0x47AD70   # accent6  (hex 70AD47)
$tcs.Item(11).RGB = 0xC16305   # hlink    (hex 0563C1)
$tcs.Item(12).RGB = 0x724F95   # folHlink (hex 954F72)
